$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in three "Razon social" (column E) entries ---
# Stray commas used as name separators were mistakenly turned into periods
# by the same pass that reformatted the floating point "Importe" values below.
$ws.Range("E98").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E178").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E206").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix "Importe" (column H) floating point formatting ---
# Values were stored as localized (Argentine) text, e.g. "11.836,00"
# ("." thousands separator, "," decimal separator). Rewrite them as plain
# floating point text, e.g. "11836.00", while keeping them as text cells
# (the column holds text, not numeric, values in this workbook).
$importeValues = @(
    "950.00",
    "11836.00",
    "8740.00",
    "930.00",
    "6100.00",
    "42195.00",
    "149790.00",
    "11450.00",
    "212270.00",
    "10500.00",
    "70000.00",
    "418000.00",
    "12293.60",
    "1820.00",
    "3540.00",
    "1695.00",
    "11820.00",
    "1170.00",
    "7988.76",
    "331800.00",
    "532946.30",
    "79412.00",
    "259359.57",
    "19535.92",
    "64418.00",
    "3019.00",
    "31920.00",
    "88008.83",
    "399.00",
    "30981.39",
    "35561.54",
    "7500.00",
    "7500.00",
    "1000.00",
    "51.10",
    "35.00",
    "500.00",
    "985.70",
    "37782.27",
    "3892.00",
    "102.35",
    "9000.00",
    "4320.00",
    "1490.16",
    "279000.00",
    "1003.60",
    "17775.56",
    "37900.00",
    "31419.03",
    "996.00",
    "446.20",
    "120.00",
    "600.00",
    "9437.15",
    "13782.28",
    "6596.00",
    "14388.68",
    "3522.62",
    "6900.00",
    "6497.20",
    "47003.40",
    "420.00",
    "23288.68",
    "1100.00",
    "1219.00",
    "5180.00",
    "150.00",
    "6395.23",
    "4267.50",
    "5639.50",
    "2244.60",
    "26220.56",
    "89400.00",
    "82175.00",
    "1500.00",
    "24800.00",
    "800.00",
    "9800.00",
    "1964.82",
    "23268.65",
    "599.50",
    "100.00",
    "171.78",
    "22200.00",
    "1652.00",
    "6400.00",
    "835.20",
    "2950.00",
    "240.00",
    "1350.00",
    "5600.00",
    "24596.00",
    "3250.00",
    "1020.00",
    "518.00",
    "2010.00",
    "1130.00",
    "550.00",
    "17875.10",
    "3473.00",
    "2514.38",
    "600235.87",
    "79840.55",
    "5.06",
    "85.21",
    "30007.34",
    "226.00",
    "521.00",
    "5235.52",
    "75.00",
    "7976.40",
    "4575.00",
    "4524.00",
    "190.30",
    "6950.00",
    "2053.00",
    "14673.00",
    "40.50",
    "545.25",
    "660.00",
    "2370.00",
    "360.00",
    "20985.86",
    "2086.00",
    "1986.00",
    "389.40",
    "7845.00",
    "16503.96",
    "1904.70",
    "352.71",
    "1200.00",
    "180.00",
    "18400.00",
    "45274.00",
    "2200.00",
    "1400.00",
    "6440.00",
    "77584.00",
    "2700.00",
    "3500.00",
    "9000.00",
    "2780.00",
    "2090.00",
    "15000.00",
    "23000.00",
    "3500.00",
    "5500.00",
    "15000.00",
    "22000.00",
    "2800.00",
    "649.00",
    "46631.40",
    "391.05",
    "8371.00",
    "18492.10",
    "22616.00",
    "1080.00",
    "10000.00",
    "6000.00",
    "14000.00",
    "6000.00",
    "4500.00",
    "3000.00",
    "3500.00",
    "3000.00",
    "4000.00",
    "4500.00",
    "3500.00",
    "4500.00",
    "4500.00",
    "5000.00",
    "10500.00",
    "4000.00",
    "6000.00",
    "4500.00",
    "6000.00",
    "4000.00",
    "5000.00",
    "4000.00",
    "2500.00",
    "20000.00",
    "4000.00",
    "12700.00",
    "3750.00",
    "4800.00",
    "3100.00",
    "3072.27",
    "381.86",
    "6084.00",
    "7526.50",
    "600.00",
    "3901.00",
    "23800.00",
    "85.00",
    "3815.00",
    "32207.00",
    "54000.00",
    "9504.00",
    "1827.82",
    "10140.00",
    "3700.00",
    "28430.00",
    "3850.00",
    "11685.06",
    "18580.00",
    "4600.00",
    "7091.14",
    "234.84",
    "82099.12",
    "880.00",
    "15442.12",
    "2960.00",
    "12138.90",
    "6600.00",
    "14820.00",
    "4128.00",
    "2000.00",
    "724134.63",
    "8500.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "60000.00",
    "60000.00",
    "60000.00",
    "30000.00",
    "3249277.28",
    "21600.00",
    "2390.00",
    "1800.00",
    "2128000.00",
    "17265025.54",
    "135500.00",
    "144500.00",
    "135500.00",
    "168500.00",
    "135500.00",
    "135500.00",
    "239000.00",
    "242450.00",
    "347000.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "239000.00",
    "342500.00",
    "239000.00",
    "135500.00",
    "244000.00",
    "135500.00",
    "135500.00",
    "140250.00",
    "135500.00",
    "633501.92",
    "397800.00",
    "15200.00",
    "14755270.22",
    "415000.00",
    "16093103.65",
    "6939.83",
    "3320.24",
    "9156.00",
    "4000.00",
    "23870.00",
    "60000.00",
    "15300.00",
    "3000.00",
    "450.00",
    "44400.00",
    "12980.00",
    "4000.00",
    "850.00",
    "4800.00",
    "5500.00",
    "47301.00",
    "6362.00"
)

$importeRange = $ws.Range("H2:H282")
$importeRange.NumberFormat = "@"
for ($i = 0; $i -lt $importeValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $importeValues[$i]
}
